$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update the header label to reflect the new cutoff date (Nov 27 -> Nov 28)
$ws.Name = "Through 2021-11-28"
$ws.Range("B1").Value = "November 2021 (through November 28)"

# Update the "November" columns (one per year) with the new day's tallies.
# Column layout: B = Nov 2021, M = Nov 2020, X = Nov 2019, AI = Nov 2018,
# AT = Nov 2017, BE = Nov 2016.

# Nov 2021 column (B)
$ws.Range("B7").Value = 4
$ws.Range("B16").Value = 5
$ws.Range("B31").Value = 6
$ws.Range("B45").Value = 5
$ws.Range("B56").Value = 2
$ws.Range("B68").Value = 6

# Nov 2020 column (M)
$ws.Range("M3").Value = 17
$ws.Range("M13").Value = 6
$ws.Range("M14").Value = 3
$ws.Range("M15").Value = 6
$ws.Range("M53").Value = 1
$ws.Range("M84").Value = 3

# Nov 2019 column (X)
$ws.Range("X13").Value = 1

# Nov 2018 column (AI)
$ws.Range("AI4").Value = 7
$ws.Range("AI24").Value = 3
$ws.Range("AI35").Value = 2

# Nov 2017 column (AT)
$ws.Range("AT5").Value = 5

# Nov 2016 column (BE)
$ws.Range("BE6").Value = 2
$ws.Range("BE20").Value = 3
$ws.Range("BE28").Value = 1
